$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 54: relabel the row from a per-port row ("Santa Cruz") into the
# Monterey-area totals row, and rename the former totals label to "Totals".
$ws.Range("A54").Value = "MONTEREY AREA TOTALS"
$ws.Range("B54").Value = "Totals"

# Column A is no longer auto (best-fit) sized; give it an explicit width.
$ws.Columns.Item(1).ColumnWidth = 11.71

# Selection moves from a single cell to the whole of column A.
$ws.Range("A1:A1048576").Select()
